$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get shuffled across rows 2-27 (Fecha, Volumen,
# Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$cols = @("D", "J", "K", "L", "M", "P")

# Mapping: new row -> source row that supplies the original values
$rowMap = @{
    2  = 14
    3  = 7
    4  = 18
    5  = 2
    6  = 8
    7  = 3
    8  = 12
    9  = 25
    10 = 23
    11 = 4
    12 = 17
    13 = 10
    14 = 13
    15 = 24
    16 = 5
    17 = 22
    18 = 11
    19 = 6
    20 = 9
    21 = 27
    22 = 15
    23 = 20
    24 = 16
    25 = 19
    26 = 26
    27 = 21
}

# Snapshot current values before overwriting anything
$original = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 27; $r++) {
        $original["$col$r"] = $ws.Range("$col$r").Value2
    }
}

foreach ($col in $cols) {
    for ($newRow = 2; $newRow -le 27; $newRow++) {
        $srcRow = $rowMap[$newRow]
        $ws.Range("$col$newRow").Value = $original["$col$srcRow"]
    }
}
